$wb = $excel.ActiveWorkbook

# --- Sheet 1 (展览) : update F column "want to go" counts ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 208
$ws1.Range("F4").Value = 379
$ws1.Range("F5").Value = 1669
$ws1.Range("F6").Value = 820
$ws1.Range("F7").Value = 711
$ws1.Range("F8").Value = 2686
$ws1.Range("F10").Value = 2075
$ws1.Range("F11").Value = 843
$ws1.Range("F12").Value = 2335
$ws1.Range("F13").Value = 726
$ws1.Range("F14").Value = 6663
$ws1.Range("F17").Value = 1264
$ws1.Range("F18").Value = 146
$ws1.Range("F19").Value = 1522
$ws1.Range("F20").Value = 1343
$ws1.Range("F23").Value = 2556
$ws1.Range("F24").Value = 1637
$ws1.Range("F25").Value = 1110
$ws1.Range("F26").Value = 1020
$ws1.Range("F27").Value = 784
$ws1.Range("F28").Value = 1120
$ws1.Range("F29").Value = 259
$ws1.Range("F30").Value = 5374
$ws1.Range("F31").Value = 289
$ws1.Range("F32").Value = 1026
$ws1.Range("F33").Value = 1269
$ws1.Range("F35").Value = 3760
$ws1.Range("F36").Value = 640
$ws1.Range("F37").Value = 1704
$ws1.Range("F38").Value = 175
$ws1.Range("F39").Value = 285
$ws1.Range("F40").Value = 967
$ws1.Range("F41").Value = 1057
$ws1.Range("F42").Value = 408
$ws1.Range("F43").Value = 1775
$ws1.Range("F44").Value = 912
$ws1.Range("F45").Value = 1051
$ws1.Range("F46").Value = 513
$ws1.Range("F47").Value = 517

# --- Sheet 2 (演出) : row 6 event removed from source, rows 7-45 shift up to 6-44 ---
$ws2 = $wb.Worksheets.Item(2)
$src = $ws2.Range("B7:I45")
$dst = $ws2.Range("B6:I44")
$src.Copy($dst)
$ws2.Rows(45).Delete()

# apply updated "want to go" counts (F column) after the shift
$ws2.Range("F6").Value = 442
$ws2.Range("F7").Value = 10
$ws2.Range("F8").Value = 496
$ws2.Range("F9").Value = 16
$ws2.Range("F10").Value = 397
$ws2.Range("F11").Value = 107
$ws2.Range("F12").Value = 142
$ws2.Range("F13").Value = 108
$ws2.Range("F14").Value = 965
$ws2.Range("F15").Value = 104
$ws2.Range("F16").Value = 3
$ws2.Range("F17").Value = 15
$ws2.Range("F18").Value = 5
$ws2.Range("F19").Value = 79
$ws2.Range("F20").Value = 604
$ws2.Range("F21").Value = 254
$ws2.Range("F22").Value = 356
$ws2.Range("F23").Value = 16
$ws2.Range("F24").Value = 172
$ws2.Range("F25").Value = 84
$ws2.Range("F26").Value = 84
$ws2.Range("F27").Value = 6
$ws2.Range("F28").Value = 2
$ws2.Range("F29").Value = 308
$ws2.Range("F30").Value = 65
$ws2.Range("F31").Value = 135
$ws2.Range("F32").Value = 16
$ws2.Range("F33").Value = 216
$ws2.Range("F34").Value = 47
$ws2.Range("F35").Value = 0
$ws2.Range("F36").Value = 112
$ws2.Range("F37").Value = 24
$ws2.Range("F38").Value = 208
$ws2.Range("F39").Value = 7
$ws2.Range("F40").Value = 1
$ws2.Range("F41").Value = 9
$ws2.Range("F42").Value = 5
$ws2.Range("F43").Value = 7
$ws2.Range("F44").Value = 5

# --- Sheet 3 (本地生活) : update F column "want to go" counts ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 3301
$ws3.Range("F5").Value = 404
$ws3.Range("F7").Value = 1470
$ws3.Range("F8").Value = 768
$ws3.Range("F9").Value = 401
$ws3.Range("F10").Value = 2831
$ws3.Range("F11").Value = 314
$ws3.Range("F12").Value = 572
$ws3.Range("F13").Value = 676
$ws3.Range("F14").Value = 1195

# --- Sheet 4 (全部类型) : update F column "want to go" counts ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 404
$ws4.Range("F3").Value = 768
$ws4.Range("F4").Value = 379
$ws4.Range("F5").Value = 401
$ws4.Range("F6").Value = 2831
$ws4.Range("F7").Value = 1669
$ws4.Range("F8").Value = 442
$ws4.Range("F9").Value = 712
$ws4.Range("F10").Value = 2686
$ws4.Range("F12").Value = 2075
$ws4.Range("F13").Value = 843
$ws4.Range("F14").Value = 2335
$ws4.Range("F16").Value = 726
$ws4.Range("F17").Value = 6664
$ws4.Range("F18").Value = 572
$ws4.Range("F20").Value = 1264
$ws4.Range("F21").Value = 676
$ws4.Range("F22").Value = 146
$ws4.Range("F23").Value = 1343
$ws4.Range("F25").Value = 2556
$ws4.Range("F26").Value = 254
$ws4.Range("F27").Value = 1638
$ws4.Range("F28").Value = 84
$ws4.Range("F29").Value = 1110
$ws4.Range("F30").Value = 1120
$ws4.Range("F31").Value = 259
$ws4.Range("F32").Value = 5374
$ws4.Range("F33").Value = 289
$ws4.Range("F34").Value = 1026
$ws4.Range("F35").Value = 1269
$ws4.Range("F36").Value = 3760
$ws4.Range("F37").Value = 640
$ws4.Range("F38").Value = 308
$ws4.Range("F39").Value = 1704
$ws4.Range("F40").Value = 175
$ws4.Range("F41").Value = 65
$ws4.Range("F42").Value = 967
$ws4.Range("F43").Value = 1775
$ws4.Range("F44").Value = 912
$ws4.Range("F45").Value = 1051
$ws4.Range("F46").Value = 513
$ws4.Range("F47").Value = 517
$ws4.Range("F48").Value = 208
$ws4.Range("F49").Value = 208
